$d = $word.ActiveDocument

# The sentence "Ansible playbooks are available in the home directory." is
# built from three runs:
#   "Ansible playbo" + "oks are available in the home" + " directory."
# It needs to become:
#   "Ansible playbo" + "ok" + "s are available in the Ansible" + " directory."
# with a "_GoBack" bookmark left right where the edit ended (just before
# " directory."), and the document's old "_GoBack" bookmark (elsewhere, in an
# empty paragraph) implicitly removed -- Word only ever keeps one "_GoBack".
$target = $d.Content
$found = $target.Find.Execute("oks are available in the home")
if (-not $found) {
    throw "Could not find the text to edit"
}

$matchStart = $target.Start
$matchEnd = $target.End

# Drop two markers *before* touching any text, so the text edit below cannot
# merge across either boundary (an edit re-coalesces same-formatting runs
# everywhere that isn't pinned by a bookmark):
#  - a throwaway marker between "ok" and "s are available in the home", which
#    is where the run needs to split;
#  - the real "_GoBack" bookmark at the end of the match, which both protects
#    the following " directory." run and is its final, correct location.
$okLen = "ok".Length
$splitPoint = $matchStart + $okLen
$splitGuard = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("ZZ_SplitGuard", $splitGuard) | Out-Null

$goBackRange = $d.Range($matchEnd, $matchEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# Replace the trailing "home" with "Ansible".
$homeLen = "home".Length
$homeRange = $d.Range($matchEnd - $homeLen, $matchEnd)
$homeRange.Text = "Ansible"

# The split guard has done its job; remove it now that no further text edits
# will happen, so the "ok" / "s are available in the Ansible" split sticks
# around without leaving a visible bookmark behind.
$d.Bookmarks("ZZ_SplitGuard").Delete()
